$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8
$ws.Cells.Item($row, 1).Value = 42654.745694444442
$ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item($row, 2).Value = $true
$ws.Cells.Item($row, 3).Value = 10057.530000000001
$ws.Cells.Item($row, 4).Value = 10053.01
$ws.Cells.Item($row, 5).Value = 75.5
$ws.Cells.Item($row, 6).Value = 75.569999999999993
$ws.Cells.Item($row, 7).Value = $false
$ws.Cells.Item($row, 7).NumberFormat = "m/d/yy h:mm"
$ws.Cells.Item($row, 8).Value = 0.09
$ws.Cells.Item($row, 9).Value = $false

$ws.Columns.AutoFit() | Out-Null
